$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the shared text "economic" -> "economical" wherever it appears.
$ws.Range("E1").Value = "economical"
$ws.Range("A5").Value = "economical"

# Move the active selection from E3 to B13.
$ws.Range("B13").Select()
